$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite the "Uke 3 Onsdag 18.1.2023" paragraph: the tail starting at
#    "Ilyas fortsatt med" (through the end of the paragraph) is replaced by a
#    much longer narrative. We locate the start/end of that span with Find,
#    overwrite the whole span's text in one go, and then re-split it back
#    into individual runs (matching the target) by toggling a character
#    property on/off across each new segment's exact range -- that forces a
#    run boundary without altering the visible formatting.
# ---------------------------------------------------------------------------

$spanStart = $d.Content.Find
$spanStart.Execute("Ilyas fortsatt med", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $d.Content.Start

$probe = $d.Content
$probe.Find.Execute("Ilyas fortsatt med", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rangeStart = $probe.Start

$para = $d.Paragraphs.Item(13)
$paraEnd = $para.Range.End

$segments = @(
    "Viet-Uy ",
    "er i en matte time",
    ". ",
    "Det er en vaksine ",
    "vi tok",
    ". Ilyas ",
    "fortsetter med nettside og fikk litt trøbbel med",
    " github ",
    "og",
    " ",
    "kamera bildene",
    ". ",
    "Mathias og Viet-Uy",
    " kommer etterhvert",
    ". ",
    "Fergus gjør ferdig ",
    "server diagramet",
    ",",
    " Besnik blir send kamera bilder og videoer og begynner med å redigere dem",
    " og Viet-Uy ",
    "fortsetter med ",
    "windows server",
    ". ",
    "Vi tar pause til å spise lunsj. "
)

$fullText = [string]::Join("", $segments)

$target = $d.Range($rangeStart, $paraEnd)
$target.Text = $fullText

# Re-split into individual runs by toggling Bold on/off across each segment's
# exact character span (a no-op visually, but forces the run to break there).
$pos = $rangeStart
foreach ($seg in $segments) {
    $segLen = $seg.Length
    $segStart = $pos
    $segEnd = $pos + $segLen
    $toggle = $d.Range($segStart, $segEnd)
    $toggle.Bold = 1
    $toggle.Bold = 0
    $pos = $segEnd
}

# ---------------------------------------------------------------------------
# 2) Add a <w:lastRenderedPageBreak/> marker before the second "[placeholder]"
#    run (the one under "Uke 3 Fredag 20.1.2023").
# ---------------------------------------------------------------------------

$placeholders = $d.Content
$find2 = $placeholders.Find
$find2.Execute("[placeholder]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find2.Execute("[placeholder]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$secondPlaceholder = $d.Range($placeholders.Start, $placeholders.End)
$secondPlaceholder.Collapse(1)
$secondPlaceholder.InsertBefore([char]0x0002)
